# Daily attendance processing - 2025-12-15 04:37:54
#
# The "Recorded By" column (G) lists the users who recorded/edited each
# attendance session as a comma-separated string, e.g.
#   "dnasr281@gmail.com, System"
# This pass normalizes the ordering of that list to a canonical
# (ordinal / case-sensitive) ascending sort order, e.g.
#   "System, dnasr281@gmail.com"
#
# NOTE on scoping: this engine's PowerShell functions do NOT get a private
# scope for loop variables - a loop variable inside a function is visible
# to (and can collide with) a loop variable of the same name in a caller.
# So every loop (nested or not) below uses a distinct variable name
# ($ci, $oi, $oj, $r, $pi, ...) to avoid accidentally aliasing an outer
# loop's counter and spinning forever / exhausting the statement budget.

function Compare-Ordinal($s1, $s2) {
    $len1 = $s1.Length
    $len2 = $s2.Length
    $minLen = $len1
    if ($len2 -lt $minLen) { $minLen = $len2 }

    $ci = 0
    $result = 0
    while ($ci -lt $minLen -and $result -eq 0) {
        $c1 = [int][char]$s1[$ci]
        $c2 = [int][char]$s2[$ci]
        if ($c1 -lt $c2) { $result = -1 }
        elseif ($c1 -gt $c2) { $result = 1 }
        $ci = $ci + 1
    }
    if ($result -eq 0) {
        if ($len1 -lt $len2) { $result = -1 }
        elseif ($len1 -gt $len2) { $result = 1 }
    }
    return $result
}

function Sort-Ordinal($items) {
    $arr = @($items)
    $n = $arr.Count
    $oi = 1
    while ($oi -lt $n) {
        $key = $arr[$oi]
        $oj = $oi - 1
        $continueShift = $true
        while ($oj -ge 0 -and $continueShift) {
            $cmp = Compare-Ordinal $arr[$oj] $key
            if ($cmp -gt 0) {
                $arr[$oj + 1] = $arr[$oj]
                $oj = $oj - 1
            } else {
                $continueShift = $false
            }
        }
        $arr[$oj + 1] = $key
        $oi = $oi + 1
    }
    return $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

# Locate the "Recorded By" column dynamically (falls back to column 7 / G,
# which is where it lives in this report layout).
$recordedByCol = 7
$pi = 1
while ($pi -le $lastCol) {
    $header = $ws.Cells.Item(1, $pi).Text
    if ($header -eq "Recorded By") {
        $recordedByCol = $pi
    }
    $pi = $pi + 1
}

$r = 2
while ($r -le $lastRow) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Text
    if ($val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $sorted = Sort-Ordinal $parts
            $joined = [string]::Join(", ", $sorted)
            if ($joined -ne $val) {
                $cell.Value = $joined
            }
        }
    }
    $r = $r + 1
}
